$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels
$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "author"

# Update column B values (rows 2-6) to 5
$ws.Range("B2").Value = 5
$ws.Range("B3").Value = 5
$ws.Range("B4").Value = 5
$ws.Range("B5").Value = 5
$ws.Range("B6").Value = 5

# Update the active selection to D7
$ws.Range("D7").Select()
